# Reorders the comma-separated "Recorded By" values in column G so that any
# "System"/"system" entries are moved to the front of the list, while the
# relative order of the remaining entries (e-mail addresses, etc.) is kept
# unchanged.
#
# Example: "backup@backdoor.com, System, system" -> "System, system, backup@backdoor.com"
#          "dnasr281@gmail.com, System"           -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    if ($text.IndexOf(",") -lt 0) {
        # Nothing to reorder - single entry.
        continue
    }

    $parts = $text.Split(",") | ForEach-Object { $_.Trim() }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Length -eq 0) {
        # No "System" entry present - nothing changes.
        continue
    }

    $newParts = $systemParts + $otherParts
    $newValue = [string]::Join(", ", $newParts)

    if ($newValue -ne $text) {
        $cell.Value = $newValue
    }
}
